$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.897799999999997
$ws.Range("B7").Value = 5.212499999999997
$ws.Range("E7").Value = 16.3806
$ws.Range("A9").Value = -21.79930000000001
$ws.Range("E10").Value = 16.6931
$ws.Range("B12").Value = 5.525299999999995
$ws.Range("E13").Value = 16.53400000000001
$ws.Range("B14").Value = 5.975599999999999
$ws.Range("D15").Value = -8.947699999999996
$ws.Range("E16").Value = 15.94020000000001
$ws.Range("A18").Value = -22.28160000000002
$ws.Range("A20").Value = -18.92629999999999
$ws.Range("E20").Value = 16.07989999999998
$ws.Range("E24").Value = 16.66160000000001
$ws.Range("B26").Value = 4.010500000000004
$ws.Range("A27").Value = -21.62869999999997
$ws.Range("B27").Value = 5.241000000000007
$ws.Range("B29").Value = 5.025999999999998
$ws.Range("D33").Value = -7.934499999999996
$ws.Range("A35").Value = -18.7718
$ws.Range("D35").Value = -9.187599999999991
$ws.Range("B37").Value = 9.102600000000002
$ws.Range("B38").Value = 4.383200000000001
$ws.Range("D38").Value = -8.978399999999995
$ws.Range("E39").Value = 16.1727
$ws.Range("D43").Value = -8.2437
$ws.Range("D44").Value = -7.2748
$ws.Range("D47").Value = -7.820600000000002
$ws.Range("E47").Value = 16.73930000000001
$ws.Range("E48").Value = 17.5515
$ws.Range("B51").Value = 5.737500000000001
$ws.Range("D51").Value = -7.658399999999998
$ws.Range("B52").Value = 5.254899999999997
$ws.Range("E52").Value = 17.2073
$ws.Range("B55").Value = 5.130199999999998
$ws.Range("E56").Value = 16.67960000000001
$ws.Range("D57").Value = -8.5425
$ws.Range("D63").Value = -7.888599999999999
$ws.Range("A69").Value = -21.65709999999999
$ws.Range("B69").Value = 5.468799999999995
$ws.Range("B70").Value = 6.910700000000002
$ws.Range("D70").Value = -7.2488
$ws.Range("A76").Value = -19.48559999999999
$ws.Range("A78").Value = -19.82919999999998
$ws.Range("B81").Value = 5.477600000000003
$ws.Range("A82").Value = -21.90200000000002
$ws.Range("A83").Value = -21.9166
$ws.Range("B83").Value = 6.145100000000002
$ws.Range("E84").Value = 16.8891
$ws.Range("D88").Value = -7.498099999999995
$ws.Range("A93").Value = -20.49409999999997
$ws.Range("D99").Value = -7.992199999999996
$ws.Range("E100").Value = 16.45480000000001
$ws.Range("E101").Value = 16.76030000000001
$ws.Range("B102").Value = 8.775300000000003
